$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-ran the averaged-intensity notebook after adding a Gaussian-Quadrature
# scheme plus three new spiral sampling schemes (Spiral-90deg-10rot-5space,
# Spiral-90deg-15rot-5space, Spiral-90deg-10rot-3space). The scheme ordering in
# the output table shifted: Gaussian-Quadrature + the new spirals now sit right
# after the "Ring Perpendicular to *" rows, and the previously-existing
# NoRotation/Rotation/HexGrid rows shifted down to make room.

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.9519747158388456
$ws.Range("D10").Value = 1.342963186989
$ws.Range("E10").Value = 0.8977606274731419
$ws.Range("F10").Value = 0.9519747158388456
$ws.Range("G10").Value = 1.176649270150119
$ws.Range("H10").Value = 0.7618749660278448
$ws.Range("I10").Value = 0.9090833387038024
$ws.Range("J10").Value = 1.342963186989
$ws.Range("K10").Value = 1.120361907231071
$ws.Range("L10").Value = 1.036168311534958
$ws.Range("M10").Value = 1.006717684197125

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9716844155471395
$ws.Range("D11").Value = 0.8499678579640696
$ws.Range("E11").Value = 1.063384785178556
$ws.Range("F11").Value = 0.9716844155471395
$ws.Range("G11").Value = 0.883128892611667
$ws.Range("H11").Value = 1.212511202825416
$ws.Range("I11").Value = 1.040384539706343
$ws.Range("J11").Value = 0.8499678579640696
$ws.Range("K11").Value = 0.9566763215713128
$ws.Range("L11").Value = 0.9641803685592262
$ws.Range("M11").Value = 1.003510282305532

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9707020074869542
$ws.Range("D12").Value = 0.8508973825589918
$ws.Range("E12").Value = 1.0634896530912
$ws.Range("F12").Value = 0.9707020074869542
$ws.Range("G12").Value = 0.8835199698801467
$ws.Range("H12").Value = 1.212863633608251
$ws.Range("I12").Value = 1.040215858736938
$ws.Range("J12").Value = 0.8508973825589918
$ws.Range("K12").Value = 0.9571935178250959
$ws.Range("L12").Value = 0.9639477626560251
$ws.Range("M12").Value = 1.003614750893747

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9715377700930535
$ws.Range("D13").Value = 0.8502458133638823
$ws.Range("E13").Value = 1.063294936125409
$ws.Range("F13").Value = 0.9715377700930535
$ws.Range("G13").Value = 0.8832899050731091
$ws.Range("H13").Value = 1.212277225737125
$ws.Range("I13").Value = 1.040316574218514
$ws.Range("J13").Value = 0.8502458133638823
$ws.Range("K13").Value = 0.9567703747446454
$ws.Range("L13").Value = 0.9641540724188494
$ws.Range("M13").Value = 1.003493704101849

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.682847999999999
$ws.Range("D14").Value = 1.878520000000002
$ws.Range("E14").Value = 0.8441439999999991
$ws.Range("F14").Value = 0.682847999999999
$ws.Range("G14").Value = 1.416948000000001
$ws.Range("H14").Value = 0.6626920000000003
$ws.Range("I14").Value = 0.7968359999999998
$ws.Range("J14").Value = 1.878520000000002
$ws.Range("K14").Value = 1.361332
$ws.Range("L14").Value = 1.02209
$ws.Range("M14").Value = 1.046998

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.53
$ws.Range("D15").Value = 2.569312500000001
$ws.Range("E15").Value = 0.68
$ws.Range("F15").Value = 0.53
$ws.Range("G15").Value = 1.78
$ws.Range("H15").Value = 0.25
$ws.Range("I15").Value = 0.6290374999999999
$ws.Range("J15").Value = 2.569312500000001
$ws.Range("K15").Value = 1.624656250000001
$ws.Range("L15").Value = 1.077328125
$ws.Range("M15").Value = 1.073058333333333

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.7238698827776012
$ws.Range("D16").Value = 1.911120490086399
$ws.Range("E16").Value = 0.8107580334080021
$ws.Range("F16").Value = 0.7238698827776012
$ws.Range("G16").Value = 1.450050220851193
$ws.Range("H16").Value = 0.5620702789632017
$ws.Range("I16").Value = 0.7842249597952017
$ws.Range("J16").Value = 1.911120490086399
$ws.Range("K16").Value = 1.3609392617472
$ws.Range("L16").Value = 1.042404572262401
$ws.Range("M16").Value = 1.040348977646933

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9986997275013679
$ws.Range("D17").Value = 0.9945532762147361
$ws.Range("E17").Value = 0.9949284922500561
$ws.Range("F17").Value = 0.9986997275013679
$ws.Range("G17").Value = 0.9978811483317608
$ws.Range("H17").Value = 0.9928038756975919
$ws.Range("I17").Value = 0.9951680834053541
$ws.Range("J17").Value = 0.9945532762147361
$ws.Range("K17").Value = 0.9947408842323961
$ws.Range("L17").Value = 0.996720305866882
$ws.Range("M17").Value = 0.9956724339001445

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.00524532515088
$ws.Range("D18").Value = 0.9528563589743533
$ws.Range("E18").Value = 1.007285047498374
$ws.Range("F18").Value = 1.00524532515088
$ws.Range("G18").Value = 0.9720768080900395
$ws.Range("H18").Value = 1.026382233025984
$ws.Range("I18").Value = 1.007056365011481
$ws.Range("J18").Value = 0.9528563589743533
$ws.Range("K18").Value = 0.9800707032363636
$ws.Range("L18").Value = 0.992658014193622
$ws.Range("M18").Value = 0.9951503562918519

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.033568193963996
$ws.Range("D19").Value = 0.8490641707221313
$ws.Range("E19").Value = 1.029746599375147
$ws.Range("F19").Value = 1.033568193963996
$ws.Range("G19").Value = 0.9086648098728471
$ws.Range("H19").Value = 1.099261416493536
$ws.Range("I19").Value = 1.035480557405926
$ws.Range("J19").Value = 0.8490641707221313
$ws.Range("K19").Value = 0.9394053850486392
$ws.Range("L19").Value = 0.9864867895063176
$ws.Range("M19").Value = 0.9926309579722642

# Rows 17-19 are brand-new rows; match the existing bold/bordered/centered
# style used for the rest of column A (style index 1 in the original file).
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "Averaged intensities table updated with spiral schemes (A1:M19)."